# TC03_Bento_MultiFilter_Diagnosis-Recurrence-TumorSize-TumorGrade
# Update the Cypher "grouped_recurrence_score" filter from "51-100" to "0-5"
# in every query cell (the per-tab "query" column B and the "StatQuery"
# column C, for the CasesTab/SamplesTab/FilesTab rows), and move the
# active selection from row 4 up to row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$queryCells = @("B2", "C2", "B3", "C3", "B4", "C4")
foreach ($addr in $queryCells) {
    $range = $ws.Range($addr)
    $text = $range.Value2
    $range.Value2 = $text -replace "51-100", "0-5"
}

$ws.Range("C2").Select()
